# Add files via upload
# Rebuild the workbook so it has 3 sheets: accept, reject, font_line
# each populated with date / time / run-id rows.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename existing Sheet1 -> accept ---------------------------
$wsAccept = $wb.Worksheets.Item(1)
$wsAccept.Name = "accept"

# --- Sheet 2: reject -------------------------------------------------------
$wsReject = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsReject.Name = "reject"

# --- Sheet 3: font_line -----------------------------------------------------
$wsFont = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsFont.Name = "font_line"

# ---------------------------------------------------------------------------
# Sheet "accept"
# ---------------------------------------------------------------------------
$wsAccept.Range("A1").Value = "d/m/y"
$wsAccept.Range("B1").Value = "time"
$wsAccept.Range("C1").Value = "run"

$wsAccept.Range("A2").Value = "06-05-2023"
$wsAccept.Range("B2").Value = "20:41:08"
$wsAccept.Range("C2").Value = "th0000-0"

$wsAccept.Range("A3").Value = "06-05-2023"
$wsAccept.Range("B3").Value = "20:41:08"
$wsAccept.Range("C3").Value = "th0000-1"

$wsAccept.Range("A4").Value = "06-05-2023"
$wsAccept.Range("B4").Value = "20:41:09"
$wsAccept.Range("C4").Value = "th0000-2"

$wsAccept.Range("A1:C1").HorizontalAlignment = -4108

$wsAccept.Columns.Item(1).ColumnWidth = 10.33203125
$wsAccept.Columns.Item(3).ColumnWidth = 21.21875

$wsAccept.Range("C6").Select()

# ---------------------------------------------------------------------------
# Sheet "reject"
# ---------------------------------------------------------------------------
$wsReject.Range("A1").Value = "d/m/y"
$wsReject.Range("B1").Value = "time"
$wsReject.Range("C1").Value = "run"

$wsReject.Range("A2").Value = "06-05-2023"
$wsReject.Range("B2").Value = "20:41:08"
$wsReject.Range("C2").Value = "th0000-0"

$wsReject.Range("A3").Value = "06-05-2023"
$wsReject.Range("B3").Value = "20:41:08"
$wsReject.Range("C3").Value = "th0000-1"

$wsReject.Range("A4").Value = "06-05-2023"
$wsReject.Range("B4").Value = "20:41:09"
$wsReject.Range("C4").Value = "th0000-2"

$wsReject.Range("A5").Value = "06-05-2023"
$wsReject.Range("B5").Value = "20:41:09"
$wsReject.Range("C5").Value = "th0000-3"

$wsReject.Range("A6").Value = "06-05-2023"
$wsReject.Range("B6").Value = "20:41:09"
$wsReject.Range("C6").Value = "th0000-4"

$wsReject.Range("A7").Value = "06-05-2024"
$wsReject.Range("B7").Value = "20:41:10"
$wsReject.Range("C7").Value = "th0000-5"

$wsReject.Range("A8").Value = "06-05-2025"
$wsReject.Range("B8").Value = "20:41:11"
$wsReject.Range("C8").Value = "th0000-6"

$wsReject.Range("A9").Value = "06-05-2026"
$wsReject.Range("B9").Value = "20:41:12"
$wsReject.Range("C9").Value = "th0000-7"

$wsReject.Range("A10").Value = "06-05-2023"
$wsReject.Range("B10").Value = "20:41:08"
$wsReject.Range("C10").Value = "th0000-8"

$wsReject.Range("A11").Value = "24-05-2023"
$wsReject.Range("B11").Value = "21:21:53"
$wsReject.Range("C11").Value = "4987072024201"

$wsReject.Range("A1:C1").HorizontalAlignment = -4108

$wsReject.Columns.Item(1).ColumnWidth = 11.6640625
$wsReject.Columns.Item(2).ColumnWidth = 13.88671875
$wsReject.Columns.Item(3).ColumnWidth = 27.77734375

$wsReject.Range("A11:C11").Select()

# ---------------------------------------------------------------------------
# Sheet "font_line"
# ---------------------------------------------------------------------------
$wsFont.Range("A1").Value = "d/m/y"
$wsFont.Range("B1").Value = "time"
$wsFont.Range("C1").Value = "run"

$wsFont.Range("A2").Value = "24-05-2023"
$wsFont.Range("B2").Value = "13:25:04"
$wsFont.Range("C2").Value = "th0000-0"

$wsFont.Range("A3").Value = "24-05-2023"
$wsFont.Range("B3").Value = "13:25:04"
$wsFont.Range("C3").Value = "th0000-1"

$wsFont.Range("A4").Value = "24-05-2023"
$wsFont.Range("B4").Value = "21:14:06"
$wsFont.Range("C4").Value = "4987072024201"

$wsFont.Range("A1:C1").HorizontalAlignment = -4108

$wsFont.Columns.Item(1).ColumnWidth = 14.44140625
$wsFont.Columns.Item(2).ColumnWidth = 19.5546875
$wsFont.Columns.Item(3).ColumnWidth = 23.33203125

$wsFont.Range("A5").Select()

# Make "accept" the active sheet/tab again, matching tabSelected=1 on sheet1
$wsAccept.Activate()
